# Apply cell-value updates to the Asura_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 42144.355
$ws.Range("I21").Value = 45000
$ws.Range("K21").Value = 45000
$ws.Range("M21").Value = -44532
# Row 23
$ws.Range("H23").Value = 42144.355
$ws.Range("I23").Value = 45000
$ws.Range("K23").Value = 45000
$ws.Range("M23").Value = -44766
# Row 92
$ws.Range("H92").Value = 5093.5
$ws.Range("I92").Value = 6641.3335
$ws.Range("J92").Value = 450
$ws.Range("K92").Value = 6641.3335
$ws.Range("L92").Value = 450
$ws.Range("M92").Value = -5393.3335
$ws.Range("N92").Value = -2946
# Row 111
$ws.Range("H111").Value = 6429.4165
$ws.Range("I111").Value = 2089.75
$ws.Range("J111").Value = 8599.25
$ws.Range("K111").Value = 6269.25
$ws.Range("L111").Value = 25797.75
$ws.Range("M111").Value = -3202.25
$ws.Range("N111").Value = -31931.75
# Row 129
$ws.Range("H129").Value = 1057.6825
$ws.Range("J129").Value = 1136.614
$ws.Range("L129").Value = 3409.842
$ws.Range("N129").Value = -13409.842

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2290.3333
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 2871
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 2871
$ws.Range("M45").Value = -1623
$ws.Range("N45").Value = -3625
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 102
$ws.Range("H102").Value = 34032.902
$ws.Range("I102").Value = 1888.3334
$ws.Range("J102").Value = 144242.86
$ws.Range("K102").Value = 1888.3334
$ws.Range("L102").Value = 144242.86
$ws.Range("M102").Value = -266.3334
$ws.Range("N102").Value = -147486.86

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 107393.734
$ws.Range("I86").Value = 2100
$ws.Range("J86").Value = 252172.62
$ws.Range("K86").Value = 2100
$ws.Range("L86").Value = 252172.62
$ws.Range("M86").Value = -977
$ws.Range("N86").Value = -254418.62
# Row 89
$ws.Range("H89").Value = 107393.734
$ws.Range("I89").Value = 2100
$ws.Range("J89").Value = 252172.62
$ws.Range("K89").Value = 10500
$ws.Range("L89").Value = 1260863.1
$ws.Range("M89").Value = -4884
$ws.Range("N89").Value = -1272095.1
# Row 94
$ws.Range("H94").Value = 53851.473
$ws.Range("I94").Value = 1141.2858
$ws.Range("J94").Value = 201440
$ws.Range("K94").Value = 1141.2858
$ws.Range("L94").Value = 201440
$ws.Range("M94").Value = -690.2858000000001
$ws.Range("N94").Value = -202342
# Row 105
$ws.Range("H105").Value = 4764045
$ws.Range("I105").Value = 5716258.5
$ws.Range("K105").Value = 5716258.5
$ws.Range("M105").Value = -5714511.5
# Row 107
$ws.Range("H107").Value = 3454.2727
$ws.Range("I107").Value = 3744.111
$ws.Range("K107").Value = 3744.111
$ws.Range("M107").Value = -1824.111

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1448.5
$ws.Range("I16").Value = 1438.2
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1438.2
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1151.2
$ws.Range("N16").Value = -2074
# Row 31
$ws.Range("H31").Value = 1826.1449
$ws.Range("I31").Value = 1351.0526
$ws.Range("J31").Value = 2408.516
$ws.Range("K31").Value = 1351.0526
$ws.Range("L31").Value = 2408.516
$ws.Range("M31").Value = -1056.0526
$ws.Range("N31").Value = -2998.516
# Row 34
$ws.Range("H34").Value = 1826.1449
$ws.Range("I34").Value = 1351.0526
$ws.Range("J34").Value = 2408.516
$ws.Range("K34").Value = 1351.0526
$ws.Range("L34").Value = 2408.516
$ws.Range("M34").Value = -1149.0526
$ws.Range("N34").Value = -2812.516
# Row 99
$ws.Range("H99").Value = 2137.5
$ws.Range("I99").Value = 2137.5
$ws.Range("K99").Value = 2137.5
$ws.Range("M99").Value = -639.5
# Row 113
$ws.Range("H113").Value = 1448.5
$ws.Range("I113").Value = 1438.2
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1438.2
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 731.8
$ws.Range("N113").Value = -5840
# Row 126
$ws.Range("H126").Value = 2137.5
$ws.Range("I126").Value = 2137.5
$ws.Range("K126").Value = 6412.5
$ws.Range("M126").Value = -3942.5

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 830.28125
$ws.Range("I68").Value = 607.2692
$ws.Range("J68").Value = 1796.6666
$ws.Range("K68").Value = 1821.8076
$ws.Range("L68").Value = 5389.9998
$ws.Range("M68").Value = -1010.8076
$ws.Range("N68").Value = -7011.9998
# Row 71
$ws.Range("H71").Value = 830.28125
$ws.Range("I71").Value = 607.2692
$ws.Range("J71").Value = 1796.6666
$ws.Range("K71").Value = 5465.422799999999
$ws.Range("L71").Value = 16169.9994
$ws.Range("M71").Value = -1409.422799999999
$ws.Range("N71").Value = -24281.9994
# Row 75
$ws.Range("H75").Value = 11621.429
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 11621.429
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 34864.287
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -36860.287
# Row 78
$ws.Range("H78").Value = 11621.429
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 11621.429
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 104592.861
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -114576.861
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
# Row 102
$ws.Range("H102").Value = 7983
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 7983
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 23949
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -28817
# Row 107
$ws.Range("H107").Value = 1257.0758
$ws.Range("I107").Value = 814.90247
$ws.Range("J107").Value = 1982.24
$ws.Range("K107").Value = 2444.70741
$ws.Range("L107").Value = 5946.72
$ws.Range("M107").Value = -524.70741
$ws.Range("N107").Value = -9786.720000000001
# Row 131
$ws.Range("H131").Value = 2621.6785
$ws.Range("I131").Value = 490.1111
$ws.Range("J131").Value = 3631.3684
$ws.Range("K131").Value = 1470.3333
$ws.Range("L131").Value = 10894.1052
$ws.Range("M131").Value = 3569.6667
$ws.Range("N131").Value = -20974.1052
# Row 140
$ws.Range("H140").Value = 2009.5834
$ws.Range("I140").Value = 836.6667
$ws.Range("J140").Value = 3182.5
$ws.Range("K140").Value = 2510.0001
$ws.Range("L140").Value = 9547.5
$ws.Range("M140").Value = 2669.9999
$ws.Range("N140").Value = -19907.5

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 5134.6816
$ws.Range("I107").Value = 7842.9287
$ws.Range("J107").Value = 395.25
$ws.Range("K107").Value = 7842.9287
$ws.Range("L107").Value = 395.25
$ws.Range("M107").Value = -5922.9287
$ws.Range("N107").Value = -4235.25
# Row 113
$ws.Range("H113").Value = 1754.7
$ws.Range("I113").Value = 1432.9166
$ws.Range("J113").Value = 2237.375
$ws.Range("K113").Value = 1432.9166
$ws.Range("L113").Value = 2237.375
$ws.Range("M113").Value = 737.0834
$ws.Range("N113").Value = -6577.375
# Row 126
$ws.Range("H126").Value = 2718.9092
$ws.Range("I126").Value = 1988.2222
$ws.Range("K126").Value = 5964.6666
$ws.Range("M126").Value = -3494.6666
# Row 132
$ws.Range("H132").Value = 2307.318
$ws.Range("I132").Value = 1728.7812
$ws.Range("J132").Value = 3850.0833
$ws.Range("K132").Value = 5186.3436
$ws.Range("L132").Value = 11550.2499
$ws.Range("M132").Value = -2656.3436
$ws.Range("N132").Value = -16610.2499

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4502
$ws.Range("I40").Value = 4004
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4004
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3868
$ws.Range("N40").Value = -5272
# Row 56
$ws.Range("H56").Value = 37698.332
$ws.Range("I56").Value = 35000
$ws.Range("J56").Value = 38469.285
$ws.Range("K56").Value = 35000
$ws.Range("L56").Value = 38469.285
$ws.Range("M56").Value = -34309
$ws.Range("N56").Value = -39851.285
# Row 122
$ws.Range("H122").Value = 10531381
$ws.Range("I122").Value = 4528
$ws.Range("J122").Value = 14290971
$ws.Range("K122").Value = 13584
$ws.Range("L122").Value = 42872913
$ws.Range("M122").Value = -11134
$ws.Range("N122").Value = -42877813
# Row 132
$ws.Range("H132").Value = 5110.4136
$ws.Range("I132").Value = 5048.2
$ws.Range("K132").Value = 15144.6
$ws.Range("M132").Value = -12614.6
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135").Value = 71380.89999999999
$ws.Range("J135").Value = 71380.89999999999
$ws.Range("L135").Value = 71380.89999999999
$ws.Range("N135").Value = -81520.89999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 58
$ws.Range("H58").Value = 37013.43
$ws.Range("I58").Value = 29500
$ws.Range("J58").Value = 40018.8
$ws.Range("K58").Value = 29500
$ws.Range("L58").Value = 40018.8
$ws.Range("M58").Value = -29192
$ws.Range("N58").Value = -40634.8
# Row 122
$ws.Range("H122").Value = 2051
$ws.Range("I122").Value = 1933
$ws.Range("K122").Value = 5799
$ws.Range("M122").Value = -3349
# Row 126
$ws.Range("H126").Value = 6986.857
$ws.Range("I126").Value = 6986.857
$ws.Range("K126").Value = 20960.571
$ws.Range("M126").Value = -18490.571
